$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking strings
# (e.g. "43.229.68", "0.998") are stored as text, matching the
# original inlineStr cell type instead of being parsed as numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '43.229.68'
$ws.Range('E2').Value = '  -0.42%  '

$ws.Range('D3').Value = '2.322.16'
$ws.Range('E3').Value = '  -2.11%  '

$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').Value = '304.76'
$ws.Range('E5').Value = '  -1.50%  '

$ws.Range('D6').Value = '100.63'
$ws.Range('E6').Value = '  -4.30%  '

$ws.Range('D7').Value = '0.511'
$ws.Range('E7').Value = '  -1.31%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').Value = '0.512'
$ws.Range('E9').Value = '  -1.15%  '

$ws.Range('D10').Value = '35.06'
$ws.Range('E10').Value = '  -3.10%  '

$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '0.0798'
$ws.Range('E11').Value = '  -1.82%  '

$ws.Range('B12').Value = 'OKB'
$ws.Range('C12').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D12').Value = '51.32'
$ws.Range('E12').Value = '  -3.86%  '

$ws.Range('E13').Value = '  +0.60%  '

$ws.Range('D14').Value = '6.82'
$ws.Range('E14').Value = '  -2.57%  '

$ws.Range('D15').Value = '2.673.69'
$ws.Range('E15').Value = '  -2.43%  '

$ws.Range('D16').Value = '15.71'
$ws.Range('E16').Value = '  +0.54%  '

$ws.Range('D17').Value = '2.273.13'
$ws.Range('E17').Value = '  -4.22%  '

$ws.Range('D18').Value = '0.804'
$ws.Range('E18').Value = '  -1.14%  '

$ws.Range('D19').Value = '43.075.42'
$ws.Range('E19').Value = '  -0.67%  '

$ws.Range('D20').Value = '11.83'
$ws.Range('E20').Value = '  -1.57%  '

$ws.Range('D21').Value = '0.0₃0905'
$ws.Range('E21').Value = '  -1.55%  '

$ws.Range('D22').Value = '6.10'
$ws.Range('E22').Value = '  -3.19%  '

$ws.Range('D23').Value = '67.57'
$ws.Range('E23').Value = '  -1.03%  '

$ws.Range('D24').Value = '238.09'
$ws.Range('E24').Value = '  -1.43%  '

$ws.Range('D25').Value = '1.98'
$ws.Range('E25').Value = '  -3.71%  '

$ws.Range('D26').Value = '2.53'
$ws.Range('E26').Value = '  -3.20%  '

$ws.Range('E27').Value = '  +0.07%  '

$ws.Range('D28').Value = '25.09'
$ws.Range('E28').Value = '  -2.88%  '

$ws.Range('D29').Value = '2.18'
$ws.Range('E29').Value = '  -5.48%  '

$ws.Range('D30').Value = '34.74'
$ws.Range('E30').Value = '  -5.63%  '

$ws.Range('E31').Value = '  +2.16%  '

$ws.Range('D32').Value = '9.23'
$ws.Range('E32').Value = '  -3.46%  '

$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.08%  '

$ws.Range('D34').Value = '5.07'
$ws.Range('E34').Value = '  -3.83%  '

$ws.Range('E35').Value = '  -4.80%  '

$ws.Range('D36').Value = '4.52'
$ws.Range('E36').Value = '  -5.19%  '

$ws.Range('B37').Value = 'Celestia'
$ws.Range('C37').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D37').Value = '16.91'
$ws.Range('E37').Value = '  -8.14%  '

$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.0707'
$ws.Range('E38').Value = '  -4.85%  '

$ws.Range('D39').Value = '2.89'
$ws.Range('E39').Value = '  -6.78%  '

$ws.Range('D40').Value = '1.82'
$ws.Range('E40').Value = '  -6.19%  '

$ws.Range('E41').Value = '  -3.81%  '

$ws.Range('E42').Value = '  -2.77%  '

$ws.Range('D43').Value = '2.41'
$ws.Range('E43').Value = '  -10.40%  '

$ws.Range('D44').Value = '1.981.87'
$ws.Range('E44').Value = '  -1.19%  '

$ws.Range('D45').Value = '0.0285'
$ws.Range('E45').Value = '  -1.82%  '

$ws.Range('D46').Value = '18.52'
$ws.Range('E46').Value = '  -6.17%  '

$ws.Range('D47').Value = '2.95'
$ws.Range('E47').Value = '  -6.31%  '

$ws.Range('D48').Value = '9.85'
$ws.Range('E48').Value = '  -7.29%  '

$ws.Range('D49').Value = '55.12'
$ws.Range('E49').Value = '  -5.01%  '

$ws.Range('D50').Value = '4.84'
$ws.Range('E50').Value = '  +2.77%  '

$ws.Range('D51').Value = '2.543.88'
$ws.Range('E51').Value = '  -0.50%  '

# Restore column D to the default (unstyled) appearance now that
# the values are committed as text, so no stray number format
# lingers on the cells themselves.
$dRange.Style = "Normal"

Write-Output "cryptos list updated"
